$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.1986156168920132
$ws.Range("D2").Value = 0.009773962508827339
$ws.Range("E2").Value = 0.1088980828063866
$ws.Range("F2").Value = 0.6704810424270136
$ws.Range("G2").Value = 0.5229827605088246
$ws.Range("H2").Value = 0.603300824350768
$ws.Range("I2").Value = 0.6564558713883173
$ws.Range("L2").Value = 0.1219626595168339
$ws.Range("M2").Value = 0.9744115844606114
$ws.Range("N2").Value = 1.833614962414231
$ws.Range("O2").Value = 2.218712525494084

# Row 3
$ws.Range("C3").Value = 0.2005261834809531
$ws.Range("D3").Value = 0.00927455425081547
$ws.Range("E3").Value = 0.1109095939905602
$ws.Range("F3").Value = 0.6518047603111938
$ws.Range("G3").Value = 0.5042366804421761
$ws.Range("H3").Value = 0.5987962019289483
$ws.Range("I3").Value = 0.6461936379585751
$ws.Range("L3").Value = 0.1246203754443238
$ws.Range("M3").Value = 0.8722858148983335
$ws.Range("N3").Value = 1.681460195919669
$ws.Range("O3").Value = 2.16918886876698

# Row 4
$ws.Range("C4").Value = 0.2018060871376619
$ws.Range("D4").Value = 0.008964929135597544
$ws.Range("E4").Value = 0.1122182846828412
$ws.Range("F4").Value = 0.6407929812252888
$ws.Range("G4").Value = 0.4931116870684775
$ws.Range("H4").Value = 0.5963692529084739
$ws.Range("I4").Value = 0.6403401554757977
$ws.Range("L4").Value = 0.1263448540712373
$ws.Range("M4").Value = 0.8093752377217953
$ws.Range("N4").Value = 1.588183730154157
$ws.Range("O4").Value = 2.140293841562737

# Row 5
$ws.Range("C5").Value = 0.2023544914193707
$ws.Range("D5").Value = 0.008838011748757424
$ws.Range("E5").Value = 0.1127700585966318
$ws.Range("F5").Value = 0.6364199651992521
$ws.Range("G5").Value = 0.4886747388048889
$ws.Range("H5").Value = 0.5954654495447045
$ws.Range("I5").Value = 0.6380671691702418
$ws.Range("L5").Value = 0.1270708319890672
$ws.Range("M5").Value = 0.7836897749842819
$ws.Range("N5").Value = 1.550214425728313
$ws.Range("O5").Value = 2.128898615365614

# Row 6
$ws.Range("C6").Value = 0.2024471734692099
$ws.Range("D6").Value = 0.008816892604706084
$ws.Range("E6").Value = 0.1128627947389775
$ws.Range("F6").Value = 0.635700732016673
$ws.Range("G6").Value = 0.487943811793258
$ws.Range("H6").Value = 0.5953205187908708
$ws.Range("I6").Value = 0.6376965219742061
$ws.Range("L6").Value = 0.1271927820663761
$ws.Range("M6").Value = 0.7794218484602169
$ws.Range("N6").Value = 1.543912301261429
$ws.Range("O6").Value = 2.127029362730951

# Row 7
$ws.Range("C7").Value = 0.20181337449484
$ws.Range("D7").Value = 0.008963220479348166
$ws.Range("E7").Value = 0.1122256513592328
$ws.Range("F7").Value = 0.6407335423413798
$ws.Range("G7").Value = 0.4930514581670025
$ws.Range("H7").Value = 0.5963567189695169
$ws.Range("I7").Value = 0.6403090465186239
$ws.Range("L7").Value = 0.1263545508414041
$ws.Range("M7").Value = 0.8090290290405022
$ws.Range("N7").Value = 1.587671488210901
$ws.Range("O7").Value = 2.14013862490313

# Row 8
$ws.Range("C8").Value = 0.199252197510873
$ws.Range("D8").Value = 0.009602391225293161
$ws.Range("E8").Value = 0.1095763498428033
$ws.Range("F8").Value = 0.6639468019719317
$ws.Range("G8").Value = 0.5164389942983405
$ws.Range("H8").Value = 0.6016772733190408
$ws.Range("I8").Value = 0.6528244408247446
$ws.Range("L8").Value = 0.1228597825555138
$ws.Range("M8").Value = 0.9392427051856913
$ws.Range("N8").Value = 1.781124150404196
$ws.Range("O8").Value = 2.201322396670832

# Row 9
$ws.Range("C9").Value = 0.1950783245289998
$ws.Range("D9").Value = 0.01083179459005379
$ws.Range("E9").Value = 0.1049668442733194
$ws.Range("F9").Value = 0.7130935386871045
$ws.Range("G9").Value = 0.5653736785633896
$ws.Range("H9").Value = 0.6148023790852477
$ws.Range("I9").Value = 0.6809291478748491
$ws.Range("L9").Value = 0.1167438203586824
$ws.Range("M9").Value = 1.192861315690038
$ws.Range("N9").Value = 2.161470496490892
$ws.Range("O9").Value = 2.333344773332698

# Row 10
$ws.Range("C10").Value = 0.1925306485703331
$ws.Range("D10").Value = 0.01172005900377471
$ws.Range("E10").Value = 0.1019395022864904
$ws.Range("F10").Value = 0.7514320853191379
$ws.Range("G10").Value = 0.6032246018390879
$ws.Range("H10").Value = 0.6260916109420975
$ws.Range("I10").Value = 0.7037673354509479
$ws.Range("L10").Value = 0.1127028494060531
$ws.Range("M10").Value = 1.378018550803233
$ws.Range("N10").Value = 2.441298695674277
$ws.Range("O10").Value = 2.437749169629228

# Row 11
$ws.Range("C11").Value = 0.191484647808899
$ws.Range("D11").Value = 0.01212083321892976
$ws.Range("E11").Value = 0.1006407548837007
$ws.Range("F11").Value = 0.7693622139460246
$ws.Range("G11").Value = 0.6208620600553729
$ws.Range("H11").Value = 0.6315861304911436
$ws.Range("I11").Value = 0.7146365727812878
$ws.Range("L11").Value = 0.1109633789319728
$ws.Range("M11").Value = 1.461971393641804
$ws.Range("N11").Value = 2.568641326600527
$ws.Range("O11").Value = 2.486869105245944

# Row 12
$ws.Range("C12").Value = 0.1911048252000853
$ws.Range("D12").Value = 0.01227211368609815
$ws.Range("E12").Value = 0.100160269166409
$ws.Range("F12").Value = 0.7762225780744529
$ws.Range("G12").Value = 0.627601514093783
$ws.Range("H12").Value = 0.6337184496340171
$ws.Range("I12").Value = 0.7188217604641238
$ws.Range("L12").Value = 0.1103189497271835
$ws.Range("M12").Value = 1.49372002156278
$ws.Range("N12").Value = 2.616865377633587
$ws.Range("O12").Value = 2.505704278029555

# Row 13
$ws.Range("C13").Value = 0.1911859024253388
$ws.Range("D13").Value = 0.01223955440562463
$ws.Range("E13").Value = 0.1002632460754231
$ws.Range("F13").Value = 0.7747419323290359
$ws.Range("G13").Value = 0.6261473552242478
$ws.Range("H13").Value = 0.6332569181310532
$ws.Range("I13").Value = 0.7179173223705106
$ws.Range("L13").Value = 0.1104571033652508
$ws.Range("M13").Value = 1.486884320512971
$ws.Range("N13").Value = 2.606479442826696
$ws.Range("O13").Value = 2.501637344761377

# Row 14
$ws.Range("C14").Value = 0.1914530732060911
$ws.Range("D14").Value = 0.01213328890761289
$ws.Range("E14").Value = 0.1006009978603686
$ws.Range("F14").Value = 0.7699252039425346
$ws.Range("G14").Value = 0.6214153041620705
$ws.Range("H14").Value = 0.6317605220129963
$ws.Range("I14").Value = 0.7149795017181617
$ws.Range("L14").Value = 0.1109100749783813
$ws.Range("M14").Value = 1.464584239870462
$ws.Range("N14").Value = 2.572608731462992
$ws.Range("O14").Value = 2.488413981996473

# Row 15
$ws.Range("C15").Value = 0.1916188435137371
$ws.Range("D15").Value = 0.0120681349126599
$ws.Range("E15").Value = 0.1008093564030981
$ws.Range("F15").Value = 0.7669840203361389
$ws.Range("G15").Value = 0.6185246780029274
$ws.Range("H15").Value = 0.6308506655683175
$ws.Range("I15").Value = 0.7131890256589628
$ws.Range("L15").Value = 0.1111893935439294
$ws.Range("M15").Value = 1.450919180472368
$ws.Range("N15").Value = 2.551862086941469
$ws.Range("O15").Value = 2.480344853244333

# Row 16
$ws.Range("C16").Value = 0.1926012828528059
$ws.Range("D16").Value = 0.01169380023691957
$ws.Range("E16").Value = 0.1020259602322593
$ws.Range("F16").Value = 0.7502701802282559
$ws.Range("G16").Value = 0.6020804068423331
$ws.Range("H16").Value = 0.6257397585196429
$ws.Range("I16").Value = 0.7030666831214774
$ws.Range("L16").Value = 0.1128185213013695
$ws.Range("M16").Value = 1.372526240379742
$ws.Range("N16").Value = 2.432977109687215
$ws.Range("O16").Value = 2.434571844533821

# Row 17
$ws.Range("C17").Value = 0.1932329303192404
$ws.Range("D17").Value = 0.01146330580409938
$ws.Range("E17").Value = 0.1027924223602001
$ws.Range("F17").Value = 0.7401423487661276
$ws.Range("G17").Value = 0.5920998765530925
$ws.Range("H17").Value = 0.6226963582382012
$ws.Range("I17").Value = 0.6969800549021627
$ws.Range("L17").Value = 0.1138432864550958
$ws.Range("M17").Value = 1.324362003563905
$ws.Range("N17").Value = 2.360053834080134
$ws.Range("O17").Value = 2.406908521690582

# Row 18
$ws.Range("C18").Value = 0.1936068666301338
$ws.Range("D18").Value = 0.01133042128007844
$ws.Range("E18").Value = 0.1032406481989501
$ws.Range("F18").Value = 0.7343631724881163
$ws.Range("G18").Value = 0.5863987491463973
$ws.Range("H18").Value = 0.6209796623674038
$ws.Range("I18").Value = 0.6935243573041845
$ws.Range("L18").Value = 0.1144419987159577
$ws.Range("M18").Value = 1.29663338501274
$ws.Range("N18").Value = 2.318115028926229
$ws.Range("O18").Value = 2.391150299822783

# Row 19
$ws.Range("C19").Value = 0.1937352996383055
$ws.Range("D19").Value = 0.011285375865004
$ws.Range("E19").Value = 0.1033936753547025
$ws.Range("F19").Value = 0.7324143519299326
$ws.Range("G19").Value = 0.5844752049559645
$ws.Range("H19").Value = 0.6204042204210083
$ws.Range("I19").Value = 0.6923620688857213
$ws.Range("L19").Value = 0.1146463067533756
$ws.Range("M19").Value = 1.287240601212332
$ws.Range("N19").Value = 2.303916225885018
$ws.Range("O19").Value = 2.385841090721215

# Row 20
$ws.Range("C20").Value = 0.1931645900253471
$ws.Range("D20").Value = 0.01148787450371813
$ws.Range("E20").Value = 0.102710067315118
$ws.Range("F20").Value = 0.7412157027543174
$ws.Range("G20").Value = 0.5931582401386066
$ws.Range("H20").Value = 0.6230168362151005
$ws.Range("I20").Value = 0.6976233103799103
$ws.Range("L20").Value = 0.1137332360435184
$ws.Range("M20").Value = 1.329491861662149
$ws.Range("N20").Value = 2.367816185155334
$ws.Range("O20").Value = 2.409837490509915

# Row 21
$ws.Range("C21").Value = 0.1913741567041498
$ws.Range("D21").Value = 0.01216451487046655
$ws.Range("E21").Value = 0.1005014843297811
$ws.Range("F21").Value = 0.7713380758827242
$ws.Range("G21").Value = 0.6228035771044347
$ws.Range("H21").Value = 0.6321986473540733
$ws.Range("I21").Value = 0.7158405298974913
$ws.Range("L21").Value = 0.1107766384192068
$ws.Range("M21").Value = 1.471135489457865
$ws.Range("N21").Value = 2.582557358607119
$ws.Range("O21").Value = 2.492291633640264

# Row 22
$ws.Range("C22").Value = 0.1902988845434592
$ws.Range("D22").Value = 0.01260391289012475
$ws.Range("E22").Value = 0.09912406127660667
$ws.Range("F22").Value = 0.7914364655733976
$ws.Range("G22").Value = 0.6425314785346075
$ws.Range("H22").Value = 0.6385006473824149
$ws.Range("I22").Value = 0.7281502390022894
$ws.Range("L22").Value = 0.1089275406857793
$ws.Range("M22").Value = 1.56345883281314
$ws.Range("N22").Value = 2.722913600717959
$ws.Range("O22").Value = 2.547547474230214

# Row 23
$ws.Range("C23").Value = 0.1908640851212553
$ws.Range("D23").Value = 0.01236965948846347
$ws.Range("E23").Value = 0.09985316207494233
$ws.Range("F23").Value = 0.7806718454922787
$ws.Range("G23").Value = 0.6319699327706587
$ws.Range("H23").Value = 0.6351095819111663
$ws.Range("I23").Value = 0.7215433057054668
$ws.Range("L23").Value = 0.1099068037397135
$ws.Range("M23").Value = 1.514207822998131
$ws.Range("N23").Value = 2.648003397965795
$ws.Range("O23").Value = 2.517931046456795

# Row 24
$ws.Range("C24").Value = 0.1931954530535052
$ws.Range("D24").Value = 0.01147676814721876
$ws.Range("E24").Value = 0.102747276433738
$ws.Range("F24").Value = 0.740730304088558
$ws.Range("G24").Value = 0.5926796393157758
$ws.Range("H24").Value = 0.6228718453416633
$ws.Range("I24").Value = 0.6973323588250437
$ws.Range("L24").Value = 0.1137829600564388
$ws.Range("M24").Value = 1.327172772079138
$ws.Range("N24").Value = 2.364306870239375
$ws.Range("O24").Value = 2.408512849624117

# Row 25
$ws.Range("C25").Value = 0.1961164606236636
$ws.Range("D25").Value = 0.01050181408954742
$ws.Range("E25").Value = 0.1061509038519226
$ws.Range("F25").Value = 0.6994078096375489
$ws.Range("G25").Value = 0.5518040055602143
$ws.Range("H25").Value = 0.6109629882139984
$ws.Range("I25").Value = 0.6729428611371731
$ws.Range("L25").Value = 0.1183191145020691
$ws.Range("M25").Value = 1.124448612931147
$ws.Range("N25").Value = 2.058493430331339
$ws.Range("O25").Value = 2.296333206437566
